# Update the two-digit division worksheet numbers to match the
# "output generated at 1c8df47" regeneration.
$d = $word.ActiveDocument

$replacements = @(
    @("19÷2=", "56÷5="),
    @("43÷7=", "20÷8="),
    @("13÷7=", "25÷8="),
    @("66÷3=", "79÷4="),
    @("34÷2=", "79÷5="),
    @("87÷8=", "75÷5="),
    @("55÷3=", "32÷7="),
    @("10÷4=", "87÷4="),
    @("55÷2=", "67÷9="),
    @("66÷9=", "79÷5="),
    @("10÷7=", "34÷8="),
    @("49÷7=", "67÷3="),
    @("35÷6=", "92÷3="),
    @("95÷2=", "87÷3="),
    @("70÷8=", "98÷9="),
    @("23÷2=", "40÷8="),
    @("71÷4=", "60÷5="),
    @("33÷5=", "91÷7="),
    @("96÷6=", "84÷5="),
    @("25÷3=", "26÷5="),
    @("85÷5=", "39÷8="),
    @("24÷2=", "89÷5="),
    @("49÷3=", "80÷9="),
    @("42÷7=", "28÷6="),
    @("12÷3=", "91÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
